$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric/date-looking string to be stored as text (shared
# string), matching values typed in as plain literal text rather than being
# auto-coerced by Excel's number/date recognizer. We temporarily format the
# cell as Text, assign the value, then clear the formatting again so the
# cell keeps the workbook's default (unstyled) look, same as its neighbours.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 10: duplicate of patient 1 / chip-1, but re-run under a new biopsy
# number and a later report date.
Set-TextValue $ws.Range("A10") "100"
Set-TextValue $ws.Range("B10") "1"
Set-TextValue $ws.Range("C10") "1234567"
$ws.Range("D10").Value = "23B000000-A1/CHIP100.1"
Set-TextValue $ws.Range("E10") "1"
Set-TextValue $ws.Range("F10") "25-may-2023"
$ws.Range("G10").Value = "Carcinoma pulmonar no microcítico"
$ws.Range("H10").Value = 15.1
$ws.Range("I10").Value = "['MYC']"
$ws.Range("J10").Value = "[40]"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = "['49.56']"
$ws.Range("M10").Value = "['EML4-ALK.E6aA20.AB374361', 'EML4-ALK.E6bA20.AB374362']"
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 6
$ws.Range("Q10").Value = 1

# Row 11: duplicate of patient 2 / chip-1, same new biopsy batch/date.
Set-TextValue $ws.Range("A11") "100"
Set-TextValue $ws.Range("B11") "2"
Set-TextValue $ws.Range("C11") "1234567"
$ws.Range("D11").Value = "23B00000-A1/CHIP100.2"
Set-TextValue $ws.Range("E11") "1"
Set-TextValue $ws.Range("F11") "25-may-2023"
$ws.Range("G11").Value = "Carcinoma pulmonar no microcítico"
$ws.Range("H11").Value = 15.1
$ws.Range("I11").Value = "['KRAS']"
$ws.Range("J11").Value = "[35]"
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = "['66.50']"
$ws.Range("M11").Value = "[]"
$ws.Range("N11").Value = 4
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 1
